$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'236.74"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'21.92"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'5.436"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.05628"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'6.474"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'3.346"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'1.076"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.7888"
$ws.Range("D9").Style = "Normal"

$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01174"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1397"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03198"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02975"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09256"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001665"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.255"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16MCDexMCB"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04753"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17CoinExTokenCET"

$ws.Range("D19").Value = "'0.006206"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'0.005091"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'0.001050"
$ws.Range("D21").Style = "Normal"

$ws.Range("D23").Value = "'3.869"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "'2.198"
$ws.Range("D24").Style = "Normal"

$ws.Range("D27").Value = "'0.0004014"
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").Value = "'0.04116"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006956"
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003503"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1036"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").Value = "'0.008915"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005441"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'0.6759"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'0.03800"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
